# Updated cryptos list on Wed Jun 19 07:27:27 UTC 2024 with GitHub Actions
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row, and for a handful of rows whose ranking changed, also
# updates the Coin name (B) and Link (C).
#
# Note: several Price values are plain decimal numbers (e.g. "604.00",
# "1.00", "0.148"). Assigning those strings straight to .Value would make
# Excel auto-convert them into numeric cells and drop the original text
# formatting (trailing zeros, etc.), whereas the source file stores them
# as plain text. Prefixing with a leading apostrophe forces Excel to keep
# them as text; resetting .Style back to "Normal" afterwards clears the
# "quote prefix" cell style that the apostrophe trick leaves behind, so
# the cell ends up with no style index at all - matching the original
# (unstyled) text cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.496.86'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '3.544.49'
$ws.Range("E3").Value = '  +2.89%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''604.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.92%  '
$ws.Range("D6").Value = '''140.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = '3.543.67'
$ws.Range("E7").Value = '  +2.88%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '''0.494'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.55%  '
$ws.Range("E10").Value = '  +2.21%  '
$ws.Range("E11").Value = '  -5.01%  '
$ws.Range("E12").Value = '  +4.55%  '
$ws.Range("D13").Value = '4.150.99'
$ws.Range("E13").Value = '  +3.19%  '
$ws.Range("E14").Value = '  +3.62%  '
$ws.Range("D15").Value = '''27.29'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.54%  '
$ws.Range("D16").Value = '3.554.91'
$ws.Range("E16").Value = '  +2.25%  '
$ws.Range("E17").Value = '  +1.62%  '
$ws.Range("D18").Value = '65.507.14'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").Value = '''10.30'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.90%  '
$ws.Range("D20").Value = '''5.95'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.95%  '
$ws.Range("D21").Value = '''14.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.59%  '
$ws.Range("D22").Value = '''396.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("E23").Value = '  +4.15%  '
$ws.Range("D24").Value = '3.694.58'
$ws.Range("E24").Value = '  +3.09%  '
$ws.Range("D25").Value = '''74.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.11%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").Value = '''0.0000117'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +10.61%  '
$ws.Range("D28").Value = '''7.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.75%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '''1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '''2.31'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.18%  '
$ws.Range("D31").Value = '''8.35'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.80%  '
$ws.Range("D32").Value = '3.562.19'
$ws.Range("E32").Value = '  +3.29%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '''1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '''0.148'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.92%  '
$ws.Range("D35").Value = '''23.81'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.46%  '
$ws.Range("E36").Value = '  +8.54%  '
$ws.Range("D37").Value = '''7.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.28%  '
$ws.Range("D38").Value = '''169.27'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.23%  '
$ws.Range("E39").Value = '  +3.35%  '
$ws.Range("D40").Value = '''5.04'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.14%  '
$ws.Range("D41").Value = '''0.0817'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.27%  '
$ws.Range("D42").Value = '''0.833'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.04%  '
$ws.Range("D43").Value = '''26.33'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +16.75%  '
$ws.Range("D44").Value = '''42.98'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.38%  '
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("E46").Value = '  +0.39%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").Value = '''1.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.82%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '''1.70'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.86%  '
$ws.Range("D49").Value = '2.472.22'
$ws.Range("E49").Value = '  +12.19%  '
$ws.Range("E50").Value = '  +4.33%  '
$ws.Range("E51").Value = '  +17.58%  '
